# Update the "Metadata" worksheet values to reflect the merged/resolved
# ValueSet-FrMedicinalProductOnly metadata.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL property (row 2): point to the new canonical HL7 France URL.
$ws.Range("B2").Value = "https://hl7.fr/fhir/fr/medication/ValueSet/FrMedicinalProductOnly"

# Date property (row 8): refreshed to the new publication date.
$ws.Range("B8").Value = "2024-12-26T10:27:36+00:00"

# Copyright property (row 14): the copyright notice text is removed, leaving
# the value cell empty (only the "Copyright" label in column A remains).
$ws.Range("B14").ClearContents()
